$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cell values per commit diff ---
$ws.Cells.Item(60, 17).Value = 0    # Q60: detect_structure 1 -> 0
$ws.Cells.Item(74, 17).Value = 0    # Q74: detect_structure 1 -> 0
$ws.Cells.Item(772, 15).Value = 1   # O772: isPivot 0 -> 1
$ws.Cells.Item(774, 18).Value = 0   # R774: backup inlineStr -> 0
$ws.Cells.Item(775, 18).Value = 0   # R775: backup inlineStr -> 0

# --- Append new weekly rows 776:797 (2024-07-01 .. 2024-11-25) ---
$newRows = @(
    @(45474, 101, 105.3000030517578, 98.70999908447266, 103.8199996948242, 103.2890625, 268555387, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 104.5, 116.9499969482422, 102.0599975585938, 113.0199966430664, 112.442008972168, 640445986, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 113.8000030517578, 118.4000015258789, 105.5, 106.2699966430664, 105.7265319824219, 348852760, 2024, 7, 15, 0, 0, 0, 29, 1, 0, 2),
    @(45495, 105.8000030517578, 108.870002746582, 96.66999816894531, 104.620002746582, 104.0849761962891, 358525468, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 105.25, 106.6500015258789, 102.3600006103516, 103.3099975585938, 102.7816696166992, 148587466, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 100, 101.5999984741211, 97.40000152587891, 97.76999664306641, 97.26999664306641, 183309131, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 96.83999633789062, 97.79000091552734, 92.5, 94.12999725341797, 93.64861297607422, 124902896, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 94.65000152587891, 98.59999847412109, 94.5, 97.08999633789062, 97.08999633789062, 104818595, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 97.59999847412109, 99.30000305175781, 93.08999633789062, 96.12000274658203, 96.12000274658203, 128861847, 2024, 8, 26, 0, 0, 0, 35, 0, 2, 2),
    @(45537, 100, 100.5, 95.5, 96.08999633789062, 96.08999633789062, 133662622, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 96, 96.5, 93.80999755859375, 94.27999877929688, 94.27999877929688, 80029910, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(45551, 94.69000244140625, 96.98000335693359, 91.38999938964844, 94.65000152587891, 94.65000152587891, 90657903, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 94.59999847412109, 96, 92.19999694824219, 95.33999633789062, 95.33999633789062, 148682323, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(45565, 95.33999633789062, 96.19000244140625, 91.05000305175781, 93.16999816894531, 93.16999816894531, 76647942, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 93.15000152587891, 93.25, 87.72000122070312, 90.91999816894531, 90.91999816894531, 117966527, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 91, 91.83999633789062, 83.01000213623047, 84.12999725341797, 84.12999725341797, 106562622, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 84.5, 84.55000305175781, 75.09999847412109, 78.5, 78.5, 191227586, 2024, 10, 21, 0, 0, 0, 43, 2, 0, 0),
    @(45593, 79, 84.37999725341797, 77.69999694824219, 84.12999725341797, 84.12999725341797, 88337127, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0),
    @(45600, 84.5, 85.58999633789062, 79.5, 82.34999847412109, 82.34999847412109, 135708959, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0),
    @(45607, 82.34999847412109, 82.58000183105469, 77.26000213623047, 78.37999725341797, 78.37999725341797, 70528890, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(45614, 78.62999725341797, 80.54000091552734, 76.41000366210938, 79.16000366210938, 79.16000366210938, 60463542, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(45621, 80.98999786376953, 84.68000030517578, 80.83000183105469, 81.44000244140625, 81.44000244140625, 118798105, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0)
)

$r = 776
foreach ($row in $newRows) {
    $c = 1
    foreach ($v in $row) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $ws.Range("A" + $r).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r = $r + 1
}
